# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de
# handback has been generated/processed and that both zh-cn and de-de
# are now in sync with en-US.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

$mdFileName = "8f43d953-10ec-4737-b24c-d3f73de2a9de.md"
$mdHyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/8ed9495efb0162249a9f0c1a780ffb7150d0ff9e/e2e/8f43d953-10ec-4737-b24c-d3f73de2a9de.md"

$zhCnTargetXlf = "8f43d953-10ec-4737-b24c-d3f73de2a9de.2c0a9d68f5e89e34d7d5e4983a16db30d1c6744e.zh-cn.xlf"
$deDeTargetXlf = "8f43d953-10ec-4737-b24c-d3f73de2a9de.2c0a9d68f5e89e34d7d5e4983a16db30d1c6744e.de-de.xlf"

$zhCnHandbackDate = "2016-08-13 09:16:26"
$deDeHandbackDate = "2016-08-13 09:16:36"

# hyperlink font formatting (to match the existing "HyperLink" cell style:
# underlined, RGB FF6495ED == BGR 15570276)
$hyperlinkColor = 15570276

function Set-ExactColumnWidth {
    param($ws, $colIndex, $targetStoredWidth)
    # This runtime stores OOXML column width as:
    #   storedWidth = (floor(ColumnWidth * 6 + 0.5) + 5) / 6
    # Solve for a ColumnWidth value that reproduces the closest achievable
    # storedWidth to $targetStoredWidth (values are quantized in 1/6ths).
    $totalPx = $targetStoredWidth * 6
    $lowPx = [Math]::Floor($totalPx)
    $highPx = [Math]::Ceiling($totalPx)
    $lowW = $lowPx / 6
    $highW = $highPx / 6
    if ([Math]::Abs($lowW - $targetStoredWidth) -le [Math]::Abs($highW - $targetStoredWidth)) {
        $bestTotalPx = $lowPx
    } else {
        $bestTotalPx = $highPx
    }
    $corePx = $bestTotalPx - 5
    $lo = ($corePx - 0.5) / 6
    $hi = ($corePx + 0.5) / 6
    $cw = $lo + (($hi - $lo) * 0.5)
    $ws.Columns.Item($colIndex).ColumnWidth = $cw
}

# ---------------------------------------------------------------------
# Overview sheet: status text for both locales (columns E and F)
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

Set-ExactColumnWidth $wsOverview 5 29.9777047293527
Set-ExactColumnWidth $wsOverview 6 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdHyperlinkUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdHyperlinkUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = $hyperlinkColor

$wsZhCn.Range("J2").Value = $zhCnTargetXlf
$wsZhCn.Range("J3").Value = $zhCnTargetXlf

$wsZhCn.Range("K2").Value = $zhCnHandbackDate
$wsZhCn.Range("K3").Value = $zhCnHandbackDate

Set-ExactColumnWidth $wsZhCn 3 29.9777047293527
Set-ExactColumnWidth $wsZhCn 9 40
Set-ExactColumnWidth $wsZhCn 10 40

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdHyperlinkUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdHyperlinkUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdFileName) | Out-Null
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = $hyperlinkColor

$wsDeDe.Range("J2").Value = $deDeTargetXlf
$wsDeDe.Range("J3").Value = $deDeTargetXlf

$wsDeDe.Range("K2").Value = $deDeHandbackDate
$wsDeDe.Range("K3").Value = $deDeHandbackDate

Set-ExactColumnWidth $wsDeDe 3 29.9777047293527
Set-ExactColumnWidth $wsDeDe 9 40
Set-ExactColumnWidth $wsDeDe 10 40

Write-Host "Handback report applied"
